$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = "Teschio Fiammeggiante "
$ws.Range("B25").Value = "Stefano Tita | Clitoriders"
$ws.Range("C25").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("D25").Value = "daniel pedrotti | iMontagna"
$ws.Range("E25").Value = "Gabriele Brentari | Shark Attack"
$ws.Range("F25").Value = "Simone Miorelli | SBARX"
